# Populate the 'orders' sheet with the full order list (30 rows).
# Existing rows 2-7 are overwritten in place and new rows 8-31 are appended,
# extending the used range from A1:E7 to A1:E31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Piyush Upadhyay", "Not Process", "Pending", 4104.09),
    @(2, "Piyush Upadhyay", "Not Process", "Pending", 4004.09),
    @(3, "Ayush Saini", "Processing", "Payment Done", 4219.42),
    @(4, "Ayush Saini", "Processing", "Payment Done", 587.7),
    @(5, "Ayush Saini", "Shipped", "Payment Done", 474.22),
    @(6, "Tushar Srivastav", "Shipped", "Payment Done", 581.72),
    @(7, "Tushar Srivastav", "Out for delivery", "Payment Done", 1601.29),
    @(8, "Tushar Srivastav", "Shipped", "Payment Done", 959.33),
    @(9, "Diksha", "Cancelled", "Payment Done", 758.73),
    @(10, "Diksha", "Not Process", "Pending", 215.45),
    @(11, "Diksha", "Processing", "Payment Done", 160.95),
    @(12, "Ritika", "Delivered", "Payment Done", 128.45),
    @(13, "Ritika", "Out for delivery", "Payment Done", 515.02),
    @(14, "Abhishek", "Cancelled", "Payment Done", 1455.93),
    @(15, "Abhishek", "Shipped", "Payment Done", 1343.86),
    @(16, "Abhishek", "Processing", "Payment Done", 1398.86),
    @(17, "Abhishek", "Delivered", "Payment Done", 2202.86),
    @(18, "Adarsh", "Cancelled", "Not Done", 379.7),
    @(19, "Adarsh", "Cancelled", "Cancelled", 1251.7),
    @(20, "Adarsh", "Shipped", "Payment Done", 845.95),
    @(21, "Adarsh", "Processing", "Payment Done", 394.93),
    @(22, "Ritika", "Delivered", "Payment Done", 420.95),
    @(23, "Ankit Upadhyay", "Out for delivery", "Payment Done", 2755.5699999999997),
    @(24, "Ankit Upadhyay", "Out for delivery", "Payment Done", 2061.35),
    @(25, "Sumit Sain", "Not Process", "Payment Done", 3558.09),
    @(26, "Sumit Sain", "Not Process", "Pending", 1495.77),
    @(27, "Sumit Sain", "Out for delivery", "Payment Done", 1563.69),
    @(28, "Amit Kumar", "Out for delivery", "Payment Done", 513.73),
    @(29, "Adarsh", "Out for delivery", "Payment Done", 918.4),
    @(30, "Prashant Sharma", "Not Process", "Payment Done", 1878.79)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
